# "added cover type and duration"
#
# Replaces the "Non Smoking / Smoke-icon-attribution / Smoking /
# Cigar-icon-attribution" block (paragraphs describing the old
# "Smoking status" cover option) with a new "Family / Family-icon-
# attribution / House / House-icon-attribution" block (describing the
# new "cover type" options), directly after the existing "Friendship"
# attribution paragraph.

$d = $word.ActiveDocument

# Locate the two boundary paragraphs of the block to remove by their
# (stable) text content instead of hard-coded indices.
$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Non Smoking*") {
        $startPara = $p
    }
    if ($p.Range.Text -like "*free-icons/cigar*") {
        $endPara = $p
    }
}

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)

$family = '<a href="https://www.flaticon.com/free-icons/family" title="family icons">Family icons created by Freepik - Flaticon</a>'
$house = '<a href="https://www.flaticon.com/free-icons/house" title="house icons">House icons created by Freepik - Flaticon</a>'

$rng.Text = "Family`r" + $family + "`r" + "House`r" + $house
